$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.681.94"
$ws.Range("E2").Value = "  -14.84%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.272.89"
$ws.Range("E3").Value = "  -21.81%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "426.53"
$ws.Range("E5").Value = "  -19.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "117.91"
$ws.Range("E6").Value = "  -17.79%  "

$ws.Range("E7").Value = "  -0.33%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.457"
$ws.Range("E8").Value = "  -17.49%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.241.74"
$ws.Range("E9").Value = "  -23.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.08"
$ws.Range("E10").Value = "  -15.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0860"
$ws.Range("E11").Value = "  -20.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.292"
$ws.Range("E12").Value = "  -18.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.120"
$ws.Range("E13").Value = "  -5.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "51.677.22"
$ws.Range("E14").Value = "  -14.81%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "18.31"
$ws.Range("E15").Value = "  -19.88%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000114"
$ws.Range("E16").Value = "  -18.99%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.268.65"
$ws.Range("E17").Value = "  -22.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.83"
$ws.Range("E18").Value = "  -23.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "291.20"
$ws.Range("E19").Value = "  -19.29%  "

$ws.Range("E20").Value = "  -0.19%  "

$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.62"
$ws.Range("E21").Value = "  -26.34%  "

$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.60"
$ws.Range("E22").Value = "  -1.66%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.12"
$ws.Range("E23").Value = "  -22.99%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "52.44"
$ws.Range("E24").Value = "  -19.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.352"
$ws.Range("E25").Value = "  -22.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.139"
$ws.Range("E26").Value = "  -22.97%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("E27").Value = "  -0.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.69"
$ws.Range("E28").Value = "  -15.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0660"
$ws.Range("E29").Value = "  -22.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "140.56"
$ws.Range("E30").Value = "  -5.91%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "16.63"
$ws.Range("E31").Value = "  -16.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.29"
$ws.Range("E32").Value = "  -23.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.56"
$ws.Range("E33").Value = "  -18.39%  "

$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.37"
$ws.Range("E34").Value = "  -22.66%  "

$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.791"
$ws.Range("E35").Value = "  -21.44%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.994"
$ws.Range("E36").Value = "  -0.30%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.969"
$ws.Range("E37").Value = "  -19.45%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "31.22"
$ws.Range("E38").Value = "  -17.68%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.16"
$ws.Range("E39").Value = "  -1.60%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.04"
$ws.Range("E40").Value = "  -18.18%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.18"
$ws.Range("E41").Value = "  -20.96%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0485"
$ws.Range("E42").Value = "  -17.39%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.838.16"
$ws.Range("E43").Value = "  -20.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.507"
$ws.Range("E44").Value = "  -21.96%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0803"
$ws.Range("E45").Value = "  -12.80%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0200"
$ws.Range("E46").Value = "  -16.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.19"
$ws.Range("E47").Value = "  -25.69%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.55"
$ws.Range("E48").Value = "  -6.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.72"
$ws.Range("E49").Value = "  -25.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "14.56"
$ws.Range("E50").Value = "  -20.73%  "

$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.29"
$ws.Range("E51").Value = "  -17.51%  "
